$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted above the current row 57, pushing
# every following record down by one row (old row 57 -> new row 58, ...,
# old row 129 -> new row 130). Insert a fresh row at 57 and populate it.
$ws.Rows.Item(57).Insert()

$ws.Range("A57").Value = 4
$ws.Range("B57").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C57").Value = 'Los Lagos'
$ws.Range("D57").Value = 44413
$ws.Range("E57").Value = 10
$ws.Range("F57").Value = 100112003
$ws.Range("G57").Value = 'Ajo'
$ws.Range("H57").Value = 'Chino'
$ws.Range("I57").Value = 'Primera'
$ws.Range("J57").Value = 150
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 15000
$ws.Range("M57").Value = 15000
$ws.Range("N57").Value = '$/caja 10 kilos'
$ws.Range("O57").Value = 'China'
$ws.Range("P57").Value = 1500
$ws.Range("Q57").Value = 10
$ws.Range("R57").Value = 'Hortaliza'
